$d = $word.ActiveDocument

# Locate the first of the four "could enter ... (optional)" list paragraphs
# (district / street name / street number / building number) by its text,
# since these four paragraphs get merged into a single "address" paragraph.
$districtIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "The user could enter district*") {
        $districtIdx = $i
        break
    }
}

$pDistrict = $d.Paragraphs.Item($districtIdx)

# Replace the district paragraph's text (keep its own paragraph mark / pPr,
# i.e. the ListParagraph style + numbering) with the merged "address" text.
$newText = "The user could enter address (optional) if address is filled it must contain from 1 up to 255 alphabetical symbols inclusive"
$pTextRange = $d.Range($pDistrict.Range.Start, $pDistrict.Range.End - 1)
$pTextRange.Text = $newText

# The following three paragraphs (street name / street number / building
# number) are now entirely superseded by the merged text above, so remove
# them completely, including their paragraph marks.
$pStreetName = $d.Paragraphs.Item($districtIdx + 1)
$pBuildingNumber = $d.Paragraphs.Item($districtIdx + 3)
$delRng = $d.Range($pStreetName.Range.Start, $pBuildingNumber.Range.End)
$delRng.Delete()

# Re-create the _GoBack bookmark that used to sit inside the (now removed)
# building-number paragraph; it belongs right after "255", immediately
# before " alphabetical symbols inclusive".
$pAddress = $d.Paragraphs.Item($districtIdx)
$markerOffset = $newText.IndexOf("255") + 3
$bookmarkPos = $pAddress.Range.Start + $markerOffset
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
